# Generate Report for Handback
# ------------------------------------------------------------
# This script mirrors a localization-status "handback" run: the
# zh-cn and de-de target files just came back in sync with en-US,
# so each language tab's "Latest Target File" / "Latest Handback
# File" / "Latest Handback DateTime" columns get populated (with a
# hyperlink on the target-file cell), the Overview + per-language
# Status columns flip from "Ready for handoff" to the handed-back
# message, and the touched columns are widened so the longer text
# fits.
# ------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$fileUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/02be4f57bd13c9d40e31b02afc08df0fd230ef66/e2e/f61d9974-5934-40fa-a413-393175814941.md"
$fileName = "f61d9974-5934-40fa-a413-393175814941.md"

# ---------------- Overview sheet ----------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$overview.Columns.Item(5).ColumnWidth = 29.15
$overview.Columns.Item(6).ColumnWidth = 29.15

# ---------------- zh-cn sheet ----------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$zhcn.Range("I2").Value = $fileName
$zhcn.Range("J2").Value = "f61d9974-5934-40fa-a413-393175814941.e9026438aee3673ccd5e54d2a39cae14eb650bf9.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-31 17:15:17"

$zhcn.Range("I3").Value = $fileName
$zhcn.Range("J3").Value = "f61d9974-5934-40fa-a413-393175814941.e9026438aee3673ccd5e54d2a39cae14eb650bf9.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-31 17:15:17"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $fileUrl, [Type]::Missing, [Type]::Missing, $fileName)
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $fileUrl, [Type]::Missing, [Type]::Missing, $fileName)

$zhcn.Columns.Item(3).ColumnWidth = 29.15
$zhcn.Columns.Item(9).ColumnWidth = 39.15
$zhcn.Columns.Item(10).ColumnWidth = 39.15

# ---------------- de-de sheet ----------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Range("I2").Value = $fileName
$dede.Range("J2").Value = "f61d9974-5934-40fa-a413-393175814941.e9026438aee3673ccd5e54d2a39cae14eb650bf9.de-de.xlf"
$dede.Range("K2").Value = "2016-08-31 17:15:30"

$dede.Range("I3").Value = $fileName
$dede.Range("J3").Value = "f61d9974-5934-40fa-a413-393175814941.e9026438aee3673ccd5e54d2a39cae14eb650bf9.de-de.xlf"
$dede.Range("K3").Value = "2016-08-31 17:15:30"

$dede.Hyperlinks.Add($dede.Range("I2"), $fileUrl, [Type]::Missing, [Type]::Missing, $fileName)
$dede.Hyperlinks.Add($dede.Range("I3"), $fileUrl, [Type]::Missing, [Type]::Missing, $fileName)

$dede.Columns.Item(3).ColumnWidth = 29.15
$dede.Columns.Item(9).ColumnWidth = 39.15
$dede.Columns.Item(10).ColumnWidth = 39.15
